# Apply the "Add files via upload" edit: extend the table from column AB to
# column AL with 10 new header columns, and fill in a handful of values in
# rows 2 and 3 (including the three new "[]" placeholder columns AC:AE).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 1), columns AC:AL -------------------------------
$ws.Range("AC1").Value = "Address"
$ws.Range("AD1").Value = "Employment"
$ws.Range("AE1").Value = "Education"
$ws.Range("AF1").Value = "Nationality"
$ws.Range("AG1").Value = "List entry to US: Date"
$ws.Range("AH1").Value = "List entry to US: i94"
$ws.Range("AI1").Value = "List entry to US: Place"
$ws.Range("AJ1").Value = "List entry to US: Visa Category"
$ws.Range("AK1").Value = "List entry to US: Status Expire"
$ws.Range("AL1").Value = "Passport Expiration Date"

# --- Row 2 updates --------------------------------------------------------
$ws.Range("A2").Value = "N/A"
$ws.Range("B2").Value = "N/A"
$ws.Range("C2").Value = "N/A"
$ws.Range("AC2").Value = "[]"
$ws.Range("AD2").Value = "[]"
$ws.Range("AE2").Value = "[]"

# --- Row 3 updates --------------------------------------------------------
$ws.Range("A3").Value = "N/A"
$ws.Range("B3").Value = "N/A"
$ws.Range("C3").Value = "N/A"
$ws.Range("M3").Value = "Yangon"
$ws.Range("N3").Value = "Myanmar"
$ws.Range("R3").Value = "Buddhism"
$ws.Range("AC3").Value = "[]"
$ws.Range("AD3").Value = "[]"
$ws.Range("AE3").Value = "[]"
$ws.Range("AF3").Value = "Burmese"

# Cells O3, P3, Q3, S3, V3, W3, Y3, Z3, AA3, AB3 were emptied out entirely in
# the source edit (their <c> elements were removed, not just blanked), so
# clear them explicitly to make sure they hold no (even empty-string) value.
$ws.Range("O3").ClearContents()
$ws.Range("P3").ClearContents()
$ws.Range("Q3").ClearContents()
$ws.Range("S3").ClearContents()
$ws.Range("V3").ClearContents()
$ws.Range("W3").ClearContents()
$ws.Range("Y3").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AA3").ClearContents()
$ws.Range("AB3").ClearContents()
